$d = $word.ActiveDocument

$d.Content.Find.Execute("Doe, Johnqwe", $true, $false, $false, $false, $false, $true, 1, $false, "Migrino, Migrino", 2)
$d.Content.Find.Execute("123 Main Street Los Angeles", $true, $false, $false, $false, $false, $true, 1, $false, "Blk 16 lot 9 Acacia homes ", 2)
$d.Content.Find.Execute("555-555-5555", $true, $false, $false, $false, $false, $true, 1, $false, "09991570900", 2)
$d.Content.Find.Execute("johndoe@example.com", $true, $false, $false, $false, $false, $true, 1, $false, "noncre123@gmail.com", 2)
$d.Content.Find.Execute("2023-12-11 09:01", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-11 12:08", 2)
$d.Content.Find.Execute("2023-12-11 17:00", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-11 7:00", 2)
$d.Content.Find.Execute("Superior Room-1", $true, $false, $false, $false, $false, $true, 1, $false, "Superior Room-1, Standard Room-5", 2)
$d.Content.Find.Execute("Kubo-3", $true, $false, $false, $false, $false, $true, 1, $false, "Umbrella-2", 2)
$d.Content.Find.Execute("3000.00", $true, $false, $false, $false, $false, $true, 1, $false, "7100.00", 2)
$d.Content.Find.Execute("1500", $true, $false, $false, $false, $false, $true, 1, $false, "3550", 2)

$d.Content.Find.Execute("2000", $true, $false, $false, $false, $false, $true, 1, $false, "6500", 2)
$d.Content.Find.Execute("1000", $true, $false, $false, $false, $false, $true, 1, $false, "600", 2)
